# "Fix last image in transfer sample"
#
# 1) The cached "today" text inside the auto-updating date field
#    (type="datetimeFigureOut") on the slide master and every slide
#    layout bumps from 5/20/20 -> 5/21/20 (PowerPoint re-stamped the
#    cached field value because the deck was saved a day later).
# 2) On the last slide only, the "S:marble1" hash/value that was
#    wrongly left over from the previous frame (110) is corrected to
#    match the rest of the deck (100).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "5/20/20") {
                $shp.TextFrame.TextRange.Text = "5/21/20"
            }
        }
    }
}

# Slide master's own date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every layout hanging off the slide master has its own cached copy.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Last slide ("last image in transfer sample"): fix the stale hash/value
# of 110 back to 100 in the three "Can" shapes that still show it.
$lastSlide = $p.Slides.Item($p.Slides.Count)

$can5 = $lastSlide.Shapes.Item("Can 5")
$null = $can5.TextFrame.TextRange.Replace("hash(110)", "hash(100)")

$can20 = $lastSlide.Shapes.Item("Can 20")
$null = $can20.TextFrame.TextRange.Replace(" 110", " 100")

$can37 = $lastSlide.Shapes.Item("Can 37")
$null = $can37.TextFrame.TextRange.Replace("hash(110)", "hash(100)")
